# Remove one term excluded from generation:
# deletes the "transcription factor binding site identification assay"
# (OBI_0002021 / "TF Binding") row from Sheet1, shifting the rows below it up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole row 9 (the excluded term); rows below shift up.
$ws.Rows.Item(9).Delete()

# Columns AK:AL beyond the real data (rows 2-8) only ever held stray empty
# placeholders; clear them out.
$ws.Range("AK2:AL8").ClearContents()

# The two former placeholder rows (now rows 9-10) were already blank in the
# source data; clear them so no stray values linger after the shift.
$ws.Range("A9:AL10").ClearContents()
